# Fixed a couple bugs after resting.
#
# Weather -> Color -> Hue lookup table lives in columns H:J of Sheet1.
# Two rows had the wrong color/hue pairing:
#   - "cloudy" (row 3) was mapped to "pink"/83, should be "green"/37.
#   - "snow"   (row 6) was mapped to "blue"/70,  should be "pink"/83.
# ("blue" is no longer used anywhere, so it drops out of the shared
# string table; "green" is a brand new string.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "green"
$ws.Range("J3").Value = 37

$ws.Range("I6").Value = "pink"
$ws.Range("J6").Value = 83

# Matches the author's last selection before saving.
$ws.Range("I7").Select()
